$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing "MEDIUM THREAT - MONITOR" formatting (fill/font as
# used on the original J8 row) before its source row is removed, so the
# escalated J5 cell can reuse that style instead of Excel inventing a
# brand-new style entry.
$ws.Range("J8").Copy()
$ws.Range("J5").PasteSpecial(-4122)

# Remove the now-obsolete trailing rows (old rows 6-11). This also
# shrinks the used range / <dimension> down to A1:K5 automatically.
$ws.Rows("6:11").Delete()

# --- Row 2 (26-FEB-26) - fare figures refreshed ---------------------------
# Force text formatting first so the "DD-MMM-YY"-shaped string isn't
# silently reinterpreted as a date serial, then restore the plain
# (general-format) cell style from an already-correct neighbour cell so
# no stray style survives on the written cell.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "26-FEB-26"
$ws.Range("K2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("D2").Value = 6005
$ws.Range("E2").Value = 6018
$ws.Range("F2").Value = -13

# --- Row 3 - date & fare figures refreshed --------------------------------
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "05-MAR-26"
$ws.Range("K3").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("D3").Value = 6005
$ws.Range("E3").Value = 6018
$ws.Range("F3").Value = -13

# --- Row 4 - date & fare figures refreshed --------------------------------
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "19-MAR-26"
$ws.Range("K4").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("D4").Value = 6005
$ws.Range("E4").Value = 6018
$ws.Range("F4").Value = -13

# --- Row 5 - date & fare figures refreshed, impact escalates to MEDIUM ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "26-MAR-26"
$ws.Range("K5").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("D5").Value = 13937
$ws.Range("E5").Value = 18828
$ws.Range("F5").Value = -4891
$ws.Range("J5").Value = "MEDIUM THREAT - MONITOR"
